$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsAVL = $wb.Worksheets.Item("AVL")

# Update the Freight lifetime value for rail on the AVL sheet (C5): 20 -> 50
$wsAVL.Range("C5").Value = 50

# Add the new note explaining the calibration choice to the About sheet (row 7)
$wsAbout.Range("B7").Value = "Freight rail locomotives are almost entirely sold off into a used locomotive market and not retired from service. We assume an average lifetime of 50."

# Move cursor/selection on AVL sheet before switching away from it
[void]$wsAVL.Range("B41").Select()

# Make the About sheet the active sheet/tab, matching the saved view state
[void]$wsAbout.Activate()
[void]$wsAbout.Range("B43").Select()
